# Auto-generated edit script applying the diff changes to cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.NumberFormat = "General"
}

Set-TextValue "D2" "331.43"
Set-TextValue "E2" "0.82%"
Set-TextValue "D3" "41.44"
Set-TextValue "E3" "3.18%"
Set-TextValue "D4" "5.751"
Set-TextValue "E4" "1.14%"
Set-TextValue "D5" "0.08125"
Set-TextValue "E5" "-0.23%"
Set-TextValue "B6" "FTXToken"
Set-TextValue "C6" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D6" "2.114"
Set-TextValue "E6" "8.71%"
Set-TextValue "B7" "KuCoinToken"
Set-TextValue "C7" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue "D7" "8.712"
Set-TextValue "E7" "-0.12%"
Set-TextValue "B8" "GateToken"
Set-TextValue "C8" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D8" "4.502"
Set-TextValue "E8" "-1.70%"
Set-TextValue "D9" "2.978"
Set-TextValue "E9" "1.18%"
Set-TextValue "D10" "0.9260"
Set-TextValue "E10" "-1.93%"
Set-TextValue "D11" "0.1281"
Set-TextValue "E11" "-0.89%"
Set-TextValue "D12" "0.1953"
Set-TextValue "E12" "-2.27%"
Set-TextValue "D13" "8.815"
Set-TextValue "E13" "14.65%"
Set-TextValue "D14" "0.09183"
Set-TextValue "E14" "-0.73%"
Set-TextValue "D15" "0.03650"
Set-TextValue "E15" "5.12%"
Set-TextValue "D16" "0.1048"
Set-TextValue "E16" "9.07%"
Set-TextValue "D17" "0.001297"
Set-TextValue "E17" "-0.84%"
Set-TextValue "D18" "0.006142"
Set-TextValue "E18" "-0.50%"
Set-TextValue "E19" "-0.07%"
Set-TextValue "D20" "0.3495"
Set-TextValue "E20" "-1.13%"
Set-TextValue "D21" "0.1376"
Set-TextValue "E21" "-2.26%"
Set-TextValue "D22" "0.2602"
Set-TextValue "E22" "6.19%"
Set-TextValue "D23" "0.04420"
Set-TextValue "E23" "-0.42%"
Set-TextValue "D24" "0.001255"
Set-TextValue "E24" "0.20%"
Set-TextValue "D25" "0.004403"
Set-TextValue "E25" "1.34%"
Set-TextValue "D26" "0.0001238"
Set-TextValue "E26" "4.11%"
Set-TextValue "D39" "0.02788"
Set-TextValue "E39" "10.49%"
Set-TextValue "D40" "0.05529"
Set-TextValue "E40" "4.55%"
Set-TextValue "D41" "0.007616"
Set-TextValue "E41" "0.19%"
Set-TextValue "D42" "0.009841"
Set-TextValue "E42" "10.43%"
Set-TextValue "E43" "-1.23%"
Set-TextValue "D44" "0.002219"
Set-TextValue "E44" "7.40%"
Set-TextValue "D45" "0.01181"
Set-TextValue "E45" "9.29%"
Set-TextValue "D46" "0.00006788"
Set-TextValue "E46" "0.44%"
Set-TextValue "D47" "0.00000000748"
Set-TextValue "E47" "-0.23%"
Set-TextValue "D48" "0.002272"
Set-TextValue "E48" "26.41%"
Set-TextValue "D49" "0.003067"
Set-TextValue "E49" "6.84%"
Set-TextValue "D50" "0.00002094"
Set-TextValue "E50" "-0.23%"
Set-TextValue "D51" "0.0001994"
Set-TextValue "E51" "-0.23%"
